$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 264, shifting existing rows 264:317 down to 265:318.
$ws.Rows.Item(264).Insert()

# Populate the new row 264 with data (copy of old row 264's unchanged fields,
# plus the new values for D, J, K, L, M, P taken from the commit).
$ws.Cells.Item(264, 1).Value = 3
$ws.Cells.Item(264, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(264, 3).Value = "Coquimbo"
$ws.Cells.Item(264, 4).Value = 44637
$ws.Cells.Item(264, 5).Value = 5
$ws.Cells.Item(264, 6).Value = 100114013
$ws.Cells.Item(264, 7).Value = "Zanahoria"
$ws.Cells.Item(264, 8).Value = "Sin especificar"
$ws.Cells.Item(264, 9).Value = "Primera"
$ws.Cells.Item(264, 10).Value = 410
$ws.Cells.Item(264, 11).Value = 7500
$ws.Cells.Item(264, 12).Value = 8000
$ws.Cells.Item(264, 13).Value = 7780
$ws.Cells.Item(264, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(264, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(264, 16).Value = 389
$ws.Cells.Item(264, 17).Value = 20
$ws.Cells.Item(264, 18).Value = "Hortaliza"

# Match the date-cell style used by column D in other rows.
$ws.Cells.Item(264, 4).NumberFormat = $ws.Cells.Item(265, 4).NumberFormat
